# Update the jxls2 sample report template: the data-row placeholders were
# renamed from upper-case bean properties (${row.ORDER_ID} etc.) to
# lower-case ones (${row.order_id} etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = '${row.order_id}'
$ws.Range("B10").Value = '${row.city_name}'
$ws.Range("C10").Value = '${row.item_name}'
$ws.Range("D10").Value = '${row.order_date}'
$ws.Range("E10").Value = '${row.volume}'

# Move the active selection to E11 (matches the saved view state).
$ws.Range("E11").Select()
